# WAT new script implementation
# Adds a new test case row (WAT28 / WAT-196) to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$newRow = 36

# Copy formatting from the row above (row 35) so the new row matches
# the existing look & feel (borders, wrap text, row height) exactly.
$ws.Rows.Item(35).Copy()
$ws.Rows.Item($newRow).PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item(35).RowHeight

$ws.Cells.Item($newRow, 1).Value = "WAT28"
$ws.Cells.Item($newRow, 3).Value = "Verify that author search results display with following fields in people card, First name Last name Alternative name(s) Number of publications Organization City/state/country (when available). Journal Details"
$ws.Cells.Item($newRow, 2).Value = "WAT-196"
$ws.Cells.Item($newRow, 4).Value = "Y"
$ws.Cells.Item($newRow, 5).Value = ""

$ws.Range("B38").Select()
